# Update "paises" (countries) workbook: daily COVID data refresh.
# - Refreshes the "Datos actualizados..." timestamp in A1.
# - Updates case-count metrics (Casos totales, Nuevos casos, Casos activos,
#   Recuperados, Casos criticos, Muertes hoy, Muertes) for the countries
#   whose figures moved.
# - Because the sheet is kept sorted by "Casos totales" (column B) descending,
#   a few countries leap-frogged their neighbours; those rows' country names
#   (column A) are swapped to keep the table sorted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed timestamp (row 1 banner)
$ws.Range("A1").Value = "Datos actualizados a 18 de Octubre de 2020 a las 16:35"

# Each entry: row, country, Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes
$updates = @(
    @(4,   "Estados Unidos",        8346163, 3498,  5432663, 2689204, 0, 14,  224296),
    @(5,   "India",                 7513080, 20353, 6614737, 784167,  0, 112, 114176),
    @(17,  "Chile",                 491760,  1757,  463943,  14182,   0, 47,  13635),
    @(21,  "Alemania",              363283,  1550,  290000,  63427,   0, 3,   9856),
    @(49,  "Portugal",              99911,   1856,  59000,   38730,   0, 19,  2181),
    @(78,  "Serbia",                36160,   214,   31536,   3848,    0, 2,   776),
    @(79,  "Jordania",              36053,   0,     6773,    28950,   0, 0,   330),
    @(80,  "Birmania",              36025,   1150,  17076,   18069,   0, 42,  880),
    @(81,  "Dinamarca",             35392,   451,   29255,   5457,    0, 1,   680),
    @(82,  "Bosnia y Herzegovina",  34112,   551,   24995,   8133,    0, 3,   984),
    @(95,  "Albania",               17055,   281,   10071,   6533,    0, 3,   451),
    @(96,  "Madagascar",            16810,   0,     16215,   357,     0, 0,   238),
    @(97,  "Noruega",               16429,   60,    11863,   4288,    0, 0,   278),
    @(123, "Malaui",                5857,    5,     4742,    934,     0, 0,   181),
    @(126, "Sri Lanka",             5536,    61,    3403,    2120,    0, 0,   13),
    @(182, "Islas Feroe",           485,     2,     472,     13,      0, 0,   0),
    @(191, "Liechtenstein",         224,     7,     132,     91,      0, 0,   1),
    @(192, "Barbados",              219,     0,     195,     17,      0, 0,   7)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 1).Value = $u[1]
    $ws.Cells.Item($row, 2).Value = $u[2]
    $ws.Cells.Item($row, 3).Value = $u[3]
    $ws.Cells.Item($row, 4).Value = $u[4]
    $ws.Cells.Item($row, 5).Value = $u[5]
    $ws.Cells.Item($row, 6).Value = $u[6]
    $ws.Cells.Item($row, 7).Value = $u[7]
    $ws.Cells.Item($row, 8).Value = $u[8]
}
